$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.545
$ws.Range("B21").Value = 9.379000000000001
$ws.Range("B23").Value = 7.398999999999999
$ws.Range("B25").Value = 6.396
$ws.Range("B53").Value = 6.315
$ws.Range("B57").Value = 5.090999999999999
$ws.Range("B59").Value = 4.435
$ws.Range("B69").Value = 5.339
$ws.Range("B79").Value = 5.488
$ws.Range("B83").Value = 5.702
$ws.Range("B93").Value = 5.608
